$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (row 1). This shifts all existing data rows up by
# one (old row 2 -> new row 1, ..., old row 20 -> new row 19).
$ws.Rows.Item(1).Delete()

# Append a new data record as the new row 20.
$ws.Cells.Item(20, 1).Value = 39
$ws.Cells.Item(20, 2).Value = "Shanice"
$ws.Cells.Item(20, 3).Value = "Mccrystal"
$ws.Cells.Item(20, 4).Value = "Female"
$ws.Cells.Item(20, 5).Value = "United States"
$ws.Cells.Item(20, 6).Value = 36
$ws.Cells.Item(20, 7).Value = "21/05/2015"
$ws.Cells.Item(20, 8).Value = 2567
